{"js": "// Apply the set of text replacements recorded in the commit diff:\n// the worksheet date and each two-digit x two-digit multiplication\n// problem/answer cell get swapped for the next day's generated set.\nconst replacements = [\n  [\"2025-04-07 Monday\", \"2025-04-08 Tuesday\"],\n  [\"38\u00d766=2508\", \"65\u00d745=2925\"],\n  [\"24\u00d712=288\", \"75\u00d765=4875\"],\n  [\"91\u00d796=8736\", \"31\u00d786=2666\"],\n  [\"91\u00d780=7280\", \"57\u00d782=4674\"],\n  [\"84\u00d761=5124\", \"98\u00d745=4410\"],\n  [\"11\u00d745=495\", \"95\u00d773=6935\"],\n  [\"91\u00d713=1183\", \"64\u00d762=3968\"],\n  [\"18\u00d733=594\", \"60\u00d798=5880\"],\n  [\"22\u00d743=946\", \"12\u00d791=1092\"],\n  [\"24\u00d788=2112\", \"91\u00d781=7371\"],\n  [\"58\u00d743=2494\", \"42\u00d737=1554\"],\n  [\"28\u00d741=1148\", \"50\u00d789=4450\"],\n  [\"33\u00d757=1881\", \"65\u00d771=4615\"],\n  [\"46\u00d727=1242\", \"75\u00d735=2625\"],\n  [\"53\u00d743=2279\", \"21\u00d748=1008\"],\n  [\"98\u00d766=6468\", \"50\u00d780=4000\"],\n  [\"41\u00d760=2460\", \"63\u00d718=1134\"],\n  [\"46\u00d778=3588\", \"58\u00d716=928\"],\n  [\"50\u00d791=4550\", \"45\u00d783=3735\"],\n  [\"45\u00d771=3195\", \"72\u00d765=4680\"],\n  [\"44\u00d728=1232\", \"47\u00d742=1974\"],\n  [\"78\u00d734=2652\", \"64\u00d724=1536\"],\n  [\"26\u00d739=1014\", \"18\u00d792=1656\"],\n  [\"79\u00d715=1185\", \"91\u00d760=5460\"],\n  [\"94\u00d775=7050\", \"26\u00d765=1690\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the set of text replacements recorded in the commit diff:\n# the worksheet date and each two-digit x two-digit multiplication\n# problem/answer cell get swapped for the next day's generated set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-07 Monday\", \"2025-04-08 Tuesday\"),\n    @(\"38\u00d766=2508\", \"65\u00d745=2925\"),\n    @(\"24\u00d712=288\", \"75\u00d765=4875\"),\n    @(\"91\u00d796=8736\", \"31\u00d786=2666\"),\n    @(\"91\u00d780=7280\", \"57\u00d782=4674\"),\n    @(\"84\u00d761=5124\", \"98\u00d745=4410\"),\n    @(\"11\u00d745=495\", \"95\u00d773=6935\"),\n    @(\"91\u00d713=1183\", \"64\u00d762=3968\"),\n    @(\"18\u00d733=594\", \"60\u00d798=5880\"),\n    @(\"22\u00d743=946\", \"12\u00d791=1092\"),\n    @(\"24\u00d788=2112\", \"91\u00d781=7371\"),\n    @(\"58\u00d743=2494\", \"42\u00d737=1554\"),\n    @(\"28\u00d741=1148\", \"50\u00d789=4450\"),\n    @(\"33\u00d757=1881\", \"65\u00d771=4615\"),\n    @(\"46\u00d727=1242\", \"75\u00d735=2625\"),\n    @(\"53\u00d743=2279\", \"21\u00d748=1008\"),\n    @(\"98\u00d766=6468\", \"50\u00d780=4000\"),\n    @(\"41\u00d760=2460\", \"63\u00d718=1134\"),\n    @(\"46\u00d778=3588\", \"58\u00d716=928\"),\n    @(\"50\u00d791=4550\", \"45\u00d783=3735\"),\n    @(\"45\u00d771=3195\", \"72\u00d765=4680\"),\n    @(\"44\u00d728=1232\", \"47\u00d742=1974\"),\n    @(\"78\u00d734=2652\", \"64\u00d724=1536\"),\n    @(\"26\u00d739=1014\", \"18\u00d792=1656\"),\n    @(\"79\u00d715=1185\", \"91\u00d760=5460\"),\n    @(\"94\u00d775=7050\", \"26\u00d765=1690\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
